# Fix heavy loadtest data corrected according to confluence Performance Test Data.xls
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("heavy")

# Input corrections on the "heavy" sheet (row 6 / row 7)
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("M6").Value = 0

$ws.Range("G7").Value = 0.25
$ws.Range("M7").Value = 0.25

# Make "heavy" the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("M8").Select()
